$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7511
$ws.Range("C3").Value = 8085
$ws.Range("K3").Value = 7767
$ws.Range("K4").Value = 1634
$ws.Range("K5").Value = 551
$ws.Range("K6").Value = 8661
$ws.Range("C7").Value = 28398
$ws.Range("K7").Value = 26124

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 87
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 475
$ws.Range("K3").Value = 513
$ws.Range("K6").Value = 572
$ws.Range("K7").Value = 1704

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 193
$ws.Range("K3").Value = 193
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 553

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 391
$ws.Range("K7").Value = 1098

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 144
$ws.Range("K7").Value = 432

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 227
$ws.Range("K7").Value = 612

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 785
$ws.Range("K8").Value = 1704
$ws.Range("K10").Value = 158
$ws.Range("K18").Value = 176
$ws.Range("K19").Value = 755
$ws.Range("K20").Value = 638
$ws.Range("K26").Value = 34
$ws.Range("K29").Value = 1440
$ws.Range("K33").Value = 1098
$ws.Range("K34").Value = 150
$ws.Range("K36").Value = 337
$ws.Range("K42").Value = 964
$ws.Range("K43").Value = 214
$ws.Range("K46").Value = 52
$ws.Range("K48").Value = 328
$ws.Range("K49").Value = 148
$ws.Range("K51").Value = 335
$ws.Range("K53").Value = 327
$ws.Range("K54").Value = 514
$ws.Range("K55").Value = 289
$ws.Range("C63").Value = 280
$ws.Range("K63").Value = 71
$ws.Range("K64").Value = 156
$ws.Range("K65").Value = 612
$ws.Range("K67").Value = 1018
$ws.Range("K78").Value = 325
$ws.Range("K79").Value = 642
$ws.Range("K83").Value = 553
$ws.Range("K84").Value = 213
$ws.Range("K85").Value = 1200
$ws.Range("K87").Value = 53
$ws.Range("K89").Value = 391
$ws.Range("K95").Value = 432
$ws.Range("K96").Value = 274
$ws.Range("K97").Value = 215
$ws.Range("K98").Value = 138
$ws.Range("C101").Value = 28398
$ws.Range("K101").Value = 26124

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 366
$ws.Range("K6").Value = 289
$ws.Range("K7").Value = 1018

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 514

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 405
$ws.Range("K5").Value = 37
$ws.Range("K6").Value = 422
$ws.Range("K7").Value = 1440

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 81
$ws.Range("K4").Value = 46
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 328

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 218
$ws.Range("K3").Value = 224
$ws.Range("K6").Value = 255
$ws.Range("K7").Value = 755

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 362
$ws.Range("K7").Value = 964

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 98
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 289

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 56
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 274

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 211
$ws.Range("K3").Value = 202
$ws.Range("K7").Value = 642

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 218
$ws.Range("K7").Value = 638

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 126
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 337

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 258
$ws.Range("K4").Value = 31
$ws.Range("K7").Value = 785

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 109
$ws.Range("K7").Value = 391

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 67
$ws.Range("K6").Value = 39

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 91
$ws.Range("K3").Value = 93
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 394
$ws.Range("K3").Value = 417
$ws.Range("K4").Value = 61
$ws.Range("K6").Value = 295
$ws.Range("K7").Value = 1200

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 53
